$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PowerPlants")
$before = $wb.Worksheets.Item("FuelsExisting")
$ws.Move($before)
